# Update gh-pages to output generated at 456a3b4
# Refresh "想去人数" (F) / "最低票价" (G) counters across the three
# event-listing sheets. A couple of rows that had sold out ("已售罄")
# now have a concrete minimum price again, so their G cell flips from
# inline text back to a number.

$wb = $excel.ActiveWorkbook

# ---- 展览 (Exhibitions) ----------------------------------------------
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value  = 25787
$ws.Range("F5").Value  = 568
$ws.Range("F6").Value  = 240
$ws.Range("F7").Value  = 572
$ws.Range("F8").Value  = 164
$ws.Range("F9").Value  = 415
$ws.Range("F11").Value = 342
$ws.Range("F12").Value = 200
$ws.Range("F13").Value = 171
$ws.Range("F14").Value = 44
$ws.Range("F15").Value = 274
$ws.Range("F16").Value = 331
$ws.Range("F18").Value = 1471
$ws.Range("F19").Value = 154
$ws.Range("F20").Value = 406

# ---- 演出 (Performances) ---------------------------------------------
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 4497
$ws.Range("G2").Value = 580
$ws.Range("F6").Value  = 74
$ws.Range("F10").Value = 418
$ws.Range("F13").Value = 9
$ws.Range("F15").Value = 26

# ---- 本地生活 (Local life) ---------------------------------------------
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value  = 4890
$ws.Range("F3").Value  = 172

# ---- 全部类型 (All types) ---------------------------------------------
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value  = 4890
$ws.Range("F5").Value  = 172
$ws.Range("F6").Value  = 25787
$ws.Range("F7").Value  = 568
$ws.Range("F8").Value  = 4497
$ws.Range("G8").Value  = 580
$ws.Range("F9").Value  = 240
$ws.Range("F11").Value = 572
$ws.Range("F14").Value = 164
$ws.Range("F15").Value = 74
$ws.Range("F16").Value = 74
$ws.Range("F20").Value = 418
$ws.Range("F21").Value = 415
$ws.Range("F24").Value = 342
$ws.Range("F25").Value = 200
$ws.Range("F26").Value = 171
$ws.Range("F27").Value = 44
$ws.Range("F29").Value = 274
$ws.Range("F30").Value = 9
$ws.Range("F32").Value = 331
$ws.Range("F34").Value = 26
$ws.Range("F35").Value = 1471
$ws.Range("F36").Value = 154
$ws.Range("F38").Value = 406
